$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.100.80'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '2.106.84'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.74%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '349.90'
$ws.Range("E5").Value = '  +4.12%  '
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5161'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4477'
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.69'
$ws.Range("E9").Value = '  -4.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08966'
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.99'
$ws.Range("E12").Value = '  +5.89%  '
$ws.Range("D13").Value = '2.093.02'
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.245'
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.18'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001149'
$ws.Range("E17").Value = '  -2.57%  '
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.85'
$ws.Range("E19").Value = '  +7.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06676'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.294'
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("D23").Value = '30.201.66'
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.90'
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.354'
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").Value = '2.357.66'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.560'
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.93'
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.75'
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.183'
$ws.Range("E31").Value = '  -2.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.647'
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.274'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.987'
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.922'
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.22'
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02587'
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06855'
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2322'
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.68'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6859'
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.257'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.27'
$ws.Range("E44").Value = '  +1.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6429'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.317'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000368'
$ws.Range("E47").Value = '  +4.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.664'
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.91'
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.225'
$ws.Range("E50").Value = '  -1.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07243'
$ws.Range("E51").Value = '  +0.80%  '
